# Insert two new rows of "Frutilla" price data at rows 287-288 of the
# "Macroferia Regional de Talca" sheet, pushing the existing rows (old
# 287..345) down to 289..347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 287 down by inserting two blank rows above it.
$ws.Rows("287:288").Insert()

# New row 287: Especial quality, 2021-11-04
$ws.Range("A287").Value = 5
$ws.Range("B287").Value = 'Macroferia Regional de Talca'
$ws.Range("C287").Value = 'Maule'
$ws.Range("D287").Value = 44504
$ws.Range("E287").Value = 7
$ws.Range("F287").Value = 'Fruta'
$ws.Range("G287").Value = 100101
$ws.Range("H287").Value = 'Berries'
$ws.Range("I287").Value = 100112025
$ws.Range("J287").Value = 'Frutilla'
$ws.Range("K287").Value = 'Sin especificar'
$ws.Range("L287").Value = 'Especial'
$ws.Range("M287").Value = 300
$ws.Range("N287").Value = 7000
$ws.Range("O287").Value = 7000
$ws.Range("P287").Value = 7000
$ws.Range("Q287").Value = '$/bandeja 7 kilos'
$ws.Range("R287").Value = 'Provincia de Melipilla'
$ws.Range("S287").Value = 1000
$ws.Range("T287").Value = 7

# New row 288: Primera quality, 2021-11-04
$ws.Range("A288").Value = 5
$ws.Range("B288").Value = 'Macroferia Regional de Talca'
$ws.Range("C288").Value = 'Maule'
$ws.Range("D288").Value = 44504
$ws.Range("E288").Value = 7
$ws.Range("F288").Value = 'Fruta'
$ws.Range("G288").Value = 100101
$ws.Range("H288").Value = 'Berries'
$ws.Range("I288").Value = 100112025
$ws.Range("J288").Value = 'Frutilla'
$ws.Range("K288").Value = 'Sin especificar'
$ws.Range("L288").Value = 'Primera'
$ws.Range("M288").Value = 180
$ws.Range("N288").Value = 5000
$ws.Range("O288").Value = 5000
$ws.Range("P288").Value = 5000
$ws.Range("Q288").Value = '$/bandeja 7 kilos'
$ws.Range("R288").Value = 'Provincia de Melipilla'
$ws.Range("S288").Value = 714
$ws.Range("T288").Value = 7
